$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/18/2023  Through  9/24/2023"

# --- Data table updates (rows 14-30) ---
# Row 14
$ws.Range("M14").Value = -40
$ws.Range("N14").Value = -93.333333333333
# Row 15
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("H15").Value = -100
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = -71.428571428571
$ws.Range("N15").Value = -89.473684210526
# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 92
$ws.Range("J16").Value = 99
$ws.Range("K16").Value = -7.070707070707
$ws.Range("L16").Value = 19.480519480519
$ws.Range("M16").Value = -44.242424242424
$ws.Range("N16").Value = -82.135922330097
# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 178
$ws.Range("J17").Value = 173
$ws.Range("K17").Value = 2.890173410404
$ws.Range("L17").Value = 28.985507246376
$ws.Range("M17").Value = 24.475524475524
$ws.Range("N17").Value = -66.666666666666
# Row 18
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 48
$ws.Range("K18").Value = -14.285714285714
$ws.Range("L18").Value = -44.186046511627
$ws.Range("M18").Value = -38.461538461538
$ws.Range("N18").Value = -91.534391534391
# Row 19
$ws.Range("C19").Value = 6
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 25
$ws.Range("I19").Value = 252
$ws.Range("J19").Value = 271
$ws.Range("K19").Value = -7.011070110701
$ws.Range("L19").Value = 10.04366812227
$ws.Range("M19").Value = 129.090909090909
$ws.Range("N19").Value = 1.612903225806
# Row 20
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -22.222222222222
$ws.Range("I20").Value = 82
$ws.Range("J20").Value = 57
$ws.Range("K20").Value = 43.859649122807
$ws.Range("L20").Value = 82.222222222222
$ws.Range("M20").Value = 110.25641025641
$ws.Range("N20").Value = -64.192139737991
# Row 21
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -35
$ws.Range("F21").Value = 83
$ws.Range("G21").Value = 66
$ws.Range("H21").Value = 25.757575757575
$ws.Range("I21").Value = 659
$ws.Range("J21").Value = 674
$ws.Range("K21").Value = -2.225519287833
$ws.Range("L21").Value = 13.425129087779
$ws.Range("M21").Value = 17.259786476868
$ws.Range("N21").Value = -69.715073529411
# Row 22
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G22").Value = 1
$ws.Range("G22").NumberFormat = "#,##0"
$ws.Range("H22").Value = 400
$ws.Range("H22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I22").Value = 23
$ws.Range("J22").Value = 18
$ws.Range("K22").Value = 27.777777777777
$ws.Range("L22").Value = 76.923076923076
$ws.Range("M22").Value = 155.555555555556
# Row 24
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -45.833333333333
$ws.Range("F24").Value = 53
$ws.Range("G24").Value = 112
$ws.Range("H24").Value = -52.678571428571
$ws.Range("I24").Value = 589
$ws.Range("J24").Value = 895
$ws.Range("K24").Value = -34.189944134078
$ws.Range("L24").Value = -4.538087520259
$ws.Range("M24").Value = 98.316498316498
# Row 25
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -14.285714285714
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 40
$ws.Range("I25").Value = 233
$ws.Range("J25").Value = 264
$ws.Range("K25").Value = -11.742424242424
$ws.Range("L25").Value = 0.8658008658
$ws.Range("M25").Value = -37.533512064343
# Row 26
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G26").Value = 2
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = -45
# Row 27
$ws.Range("D27").Value = 1
$ws.Range("J27").Value = 39
$ws.Range("K27").Value = -38.461538461538
$ws.Range("L27").Value = -7.692307692307
# Row 28
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J28").Value = 13
$ws.Range("K28").Value = -53.846153846153
$ws.Range("N28").Value = -90.625
# Row 29
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J29").Value = 13
$ws.Range("K29").Value = -53.846153846153
$ws.Range("N29").Value = -90.322580645161
# Row 30
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("H30").Value = -100
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = -20
